$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.707.90'
$ws.Range("E2").Value = '  +2.53%  '

$ws.Range("D3").Value = '1.888.89'
$ws.Range("E3").Value = '  +0.70%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").Value = "'247.60"
$ws.Range("E5").Value = '  +2.32%  '

$ws.Range("E6").Value = '  +0.12%  '

$ws.Range("D7").Value = "'0.4939"
$ws.Range("E7").Value = '  +0.35%  '

$ws.Range("D8").Value = "'0.2956"
$ws.Range("E8").Value = '  +1.58%  '

$ws.Range("D9").Value = "'0.06811"
$ws.Range("E9").Value = '  +2.88%  '

$ws.Range("D10").Value = '1.887.23'
$ws.Range("E10").Value = '  +0.49%  '

$ws.Range("D11").Value = "'17.19"
$ws.Range("E11").Value = '  +2.93%  '

$ws.Range("D12").Value = "'0.07236"
$ws.Range("E12").Value = '  -0.04%  '

$ws.Range("D13").Value = "'91.74"
$ws.Range("E13").Value = '  +6.48%  '

$ws.Range("D14").Value = "'5.073"

$ws.Range("E15").Value = '  +2.07%  '

$ws.Range("D16").Value = '30.671.20'
$ws.Range("E16").Value = '  +2.42%  '

$ws.Range("D17").Value = "'0.000007984"
$ws.Range("E17").Value = '  +1.82%  '

$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = '  +0.19%  '

$ws.Range("D19").Value = "'13.22"
$ws.Range("E19").Value = '  +3.78%  '

$ws.Range("D20").Value = '2.134.51'
$ws.Range("E20").Value = '  +0.55%  '

$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = '  +0.23%  '

$ws.Range("D22").Value = "'4.832"
$ws.Range("E22").Value = '  +1.42%  '

$ws.Range("D23").Value = "'188.22"
$ws.Range("E23").Value = '  +33.68%  '

$ws.Range("D24").Value = "'6.057"
$ws.Range("E24").Value = '  +5.45%  '

$ws.Range("D25").Value = "'9.354"
$ws.Range("E25").Value = '  +3.41%  '

$ws.Range("D26").Value = "'156.29"
$ws.Range("E26").Value = '  +4.16%  '

$ws.Range("D27").Value = "'19.16"
$ws.Range("E27").Value = '  +12.85%  '

$ws.Range("D28").Value = "'1.908"
$ws.Range("E28").Value = '  -0.31%  '

$ws.Range("D29").Value = "'1.402"
$ws.Range("E29").Value = '  +0.60%  '

$ws.Range("D30").Value = "'4.296"
$ws.Range("E30").Value = '  +2.71%  '

$ws.Range("D31").Value = "'0.09007"
$ws.Range("E31").Value = '  +3.33%  '

$ws.Range("D32").Value = "'4.012"
$ws.Range("E32").Value = '  +1.50%  '

$ws.Range("D33").Value = "'0.05188"
$ws.Range("E33").Value = '  +2.81%  '

$ws.Range("D34").Value = "'0.7430"
$ws.Range("E34").Value = '  +4.63%  '

$ws.Range("D35").Value = "'1.115"
$ws.Range("E35").Value = '  +0.26%  '

$ws.Range("D36").Value = "'2.734"
$ws.Range("E36").Value = '  +2.41%  '

$ws.Range("D37").Value = "'0.01841"
$ws.Range("E37").Value = '  +3.36%  '

$ws.Range("D38").Value = "'2.666"
$ws.Range("E38").Value = '  -0.70%  '

$ws.Range("D39").Value = "'2.154"
$ws.Range("E39").Value = '  -0.56%  '

$ws.Range("D40").Value = "'0.9403"
$ws.Range("E40").Value = '  +1.29%  '

$ws.Range("D41").Value = "'0.4426"
$ws.Range("E41").Value = '  +4.56%  '

$ws.Range("D42").Value = "'105.30"
$ws.Range("E42").Value = '  +2.65%  '

$ws.Range("E43").Value = '  +0.29%  '

$ws.Range("D44").Value = "'5.771"
$ws.Range("E44").Value = '  +0.28%  '

$ws.Range("D45").Value = "'7.632"
$ws.Range("E45").Value = '  +2.93%  '

$ws.Range("D46").Value = "'0.1343"
$ws.Range("E46").Value = '  +6.15%  '

$ws.Range("D47").Value = "'0.05841"

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = "'1.428"
$ws.Range("E48").Value = '  +7.26%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = "'8.675"
$ws.Range("E49").Value = '  +5.11%  '

$ws.Range("D50").Value = "'0.3944"
$ws.Range("E50").Value = '  +4.56%  '

$ws.Range("D51").Value = "'33.51"
$ws.Range("E51").Value = '  +3.20%  '
